# Fix formatting issues introduced when scraping data:
# 1) A handful of "Razon social" text fields had commas that should have
#    been periods (typos/normalization artifacts from the scraper).
# 2) The "Importe" (amount) column was scraped with Spanish/Argentine
#    numeric formatting (dot as thousands separator, comma as decimal
#    separator) but stored as text; fix it to use a plain dot-decimal
#    textual number (no thousands separator).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix specific "Razon social" (column E) text values ---------------
$ws.Range("E106").Value = "BOFFELLI. MARIA INES"
$ws.Range("E114").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E153").Value = "RICCOTTI. MARIANA EDITH"
$ws.Range("E213").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Range("E226").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"

# --- 2) Fix numeric formatting of "Importe" (column H) amounts -----------
# Original values look like "1.234,56" (dot = thousands sep, comma =
# decimal sep). Convert to "1234.56" (dot decimal, no thousands sep),
# keeping the values as literal text (they were stored as text before,
# not as real numbers), so we force a leading quote to keep Excel from
# re-interpreting the cleaned-up digit string as a genuine number.
for ($r = 2; $r -le 294; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $old = $cell.Value2
    $new = $old.Replace(".", "").Replace(",", ".")
    $cell.Formula = "'" + $new
}
